$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("Q4").Value = 1
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 1
$ws.Range("Q6").Value = 1
$ws.Range("Q7").Value = 1
$ws.Range("Q8").Value = 1
$ws.Range("Q9").Value = 1
$ws.Range("P10").Value = 0
$ws.Range("Q10").Value = 1
$ws.Range("Q11").Value = 1
$ws.Range("P12").Value = 0
$ws.Range("Q12").Value = 1
$ws.Range("Q13").Value = 1
$ws.Range("Q14").Value = 2
$ws.Range("Q15").Value = 3
$ws.Range("Q16").Value = 3
$ws.Range("Q17").Value = 4
$ws.Range("Q18").Value = 4
$ws.Range("P19").Value = 0
$ws.Range("Q19").Value = 4
$ws.Range("Q20").Value = 4
$ws.Range("Q21").Value = 5
$ws.Range("Q22").Value = 5
$ws.Range("Q23").Value = 6
$ws.Range("O24").Value = '[''Portugal'', ''Northern Ireland'']'
$ws.Range("Q24").Value = 7
$ws.Range("Q25").Value = 8
$ws.Range("Q26").Value = 9
$ws.Range("Q27").Value = 10
$ws.Range("Q28").Value = 10
$ws.Range("Q29").Value = 11
$ws.Range("Q30").Value = 12
$ws.Range("Q31").Value = 13
$ws.Range("Q32").Value = 14
$ws.Range("Q33").Value = 14
$ws.Range("Q34").Value = 15
$ws.Range("P42").Value = 0
$ws.Range("Q42").Value = 2
$ws.Range("Q43").Value = 2
$ws.Range("Q44").Value = 3
$ws.Range("Q45").Value = 3
$ws.Range("Q46").Value = 3
$ws.Range("Q47").Value = 3
$ws.Range("P48").Value = 0
$ws.Range("Q48").Value = 3
$ws.Range("P49").Value = 0
$ws.Range("Q49").Value = 3
$ws.Range("Q50").Value = 4
$ws.Range("P51").Value = 0
$ws.Range("Q51").Value = 4
$ws.Range("P52").Value = 0
$ws.Range("Q52").Value = 4
$ws.Range("Q53").Value = 5
$ws.Range("Q54").Value = 6
$ws.Range("Q55").Value = 6
$ws.Range("Q56").Value = 6
$ws.Range("Q57").Value = 7
$ws.Range("Q58").Value = 8
$ws.Range("Q59").Value = 8
$ws.Range("Q60").Value = 9
$ws.Range("Q61").Value = 10
$ws.Range("Q62").Value = 11
$ws.Range("Q63").Value = 12
$ws.Range("Q64").Value = 13
$ws.Range("P65").Value = 0
$ws.Range("Q65").Value = 13
$ws.Range("Q66").Value = 13
$ws.Range("Q67").Value = 14
$ws.Range("Q68").Value = 14
$ws.Range("Q69").Value = 14
$ws.Range("Q70").Value = 14
$ws.Range("Q71").Value = 14
$ws.Range("Q72").Value = 14
$ws.Range("Q73").Value = 14
$ws.Range("Q74").Value = 14
$ws.Range("Q75").Value = 14
$ws.Range("P76").Value = 0
$ws.Range("Q76").Value = 14
$ws.Range("Q77").Value = 15
$ws.Range("Q78").Value = 15
$ws.Range("Q79").Value = 15
$ws.Range("Q80").Value = 16
$ws.Range("Q81").Value = 17
$ws.Range("Q82").Value = 18
$ws.Range("Q83").Value = 19
$ws.Range("Q84").Value = 20
$ws.Range("P86").Value = 0
$ws.Range("Q86").Value = 0
$ws.Range("Q87").Value = 0
$ws.Range("Q88").Value = 0
$ws.Range("Q89").Value = 1
$ws.Range("Q90").Value = 2
$ws.Range("Q91").Value = 3
$ws.Range("Q92").Value = 3
$ws.Range("Q93").Value = 4
$ws.Range("Q94").Value = 5
$ws.Range("Q95").Value = 5
$ws.Range("Q96").Value = 6
$ws.Range("Q97").Value = 7
$ws.Range("Q98").Value = 7
$ws.Range("Q99").Value = 8
$ws.Range("Q100").Value = 9
$ws.Range("Q101").Value = 9
$ws.Range("Q102").Value = 10
$ws.Range("P103").Value = 0
$ws.Range("Q103").Value = 10
$ws.Range("O104").Value = '[''Slovakia'', ''Slovenia'']'
$ws.Range("Q104").Value = 10
$ws.Range("P105").Value = 0
$ws.Range("Q105").Value = 10
$ws.Range("P106").Value = 0
$ws.Range("Q106").Value = 10
$ws.Range("Q107").Value = 11
$ws.Range("Q108").Value = 12
$ws.Range("Q109").Value = 13
$ws.Range("Q110").Value = 14
$ws.Range("Q111").Value = 15
$ws.Range("Q112").Value = 15
$ws.Range("O113").Value = '[''Netherlands'', ''Georgia'']'
$ws.Range("P113").Value = 0
$ws.Range("Q113").Value = 15
$ws.Range("Q114").Value = 16
$ws.Range("O115").Value = '[''Netherlands'', ''Georgia'']'
$ws.Range("Q115").Value = 17
